$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Company title cell from "Company1" to "CTS"
$ws.Range("E2").Value = "CTS"

# Match the saved selection state (active cell moved to E2)
$ws.Range("E2").Select()
